# Refresh the cryptos price/volume snapshot (D/E columns) with the
# latest values, as produced by the scheduled GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.355.84"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "1.839.13"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("D4").Value = "'1.016"
$ws.Range("E4").Value = "  +1.51%  "
$ws.Range("D5").Value = "'315.06"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").Value = "'0.4736"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("D8").Value = "'0.3699"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "'0.8852"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").Value = "'20.49"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "1.853.80"
$ws.Range("E12").Value = "  +2.48%  "
$ws.Range("D13").Value = "'0.07376"
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").Value = "'5.478"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "'93.24"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "'6.579"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "'0.000008846"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("D20").Value = "'14.83"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").Value = "27.388.44"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").Value = "'5.353"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "'10.72"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").Value = "2.075.71"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "'1.913"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "'152.09"
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").Value = "'2.180"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "'5.276"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").Value = "'117.92"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").Value = "'0.08933"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "'0.7620"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").Value = "'1.180"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").Value = "'4.561"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("D35").Value = "'2.944"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D38").Value = "'0.05370"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "'3.001"
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("D41").Value = "'7.294"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").Value = "'0.5357"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").Value = "'2.378"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "'0.1667"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "'8.553"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").Value = "'0.4986"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").Value = "'10.57"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").Value = "'105.25"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").Value = "'1.678"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "'0.06328"
$ws.Range("E51").Value = "  +0.71%  "
